$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.998.91"
$ws.Range("E2").Value = "  +0.33%  "

$ws.Range("D3").Value = "1.560.41"
$ws.Range("E3").Value = "  +0.55%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.33%  "

$ws.Range("E6").Value = "  +1.02%  "

$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.13"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.95%  "

$ws.Range("E9").Value = "  +0.14%  "

$ws.Range("E10").Value = "  +1.78%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0860"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.25%  "

$ws.Range("D12").Value = "1.783.14"
$ws.Range("E12").Value = "  +0.57%  "

$ws.Range("D13").Value = "1.536.88"
$ws.Range("E13").Value = "  -1.03%  "

$ws.Range("E14").Value = "  +0.97%  "

$ws.Range("E15").Value = "  +1.32%  "

$ws.Range("E16").Value = "  +0.71%  "

$ws.Range("D17").Value = "26.999.26"
$ws.Range("E17").Value = "  +0.44%  "

$ws.Range("E18").Value = "  +2.30%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "217.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.17%  "

$ws.Range("E20").Value = "  +2.11%  "

$ws.Range("E21").Value = "  -0.13%  "

$ws.Range("E22").Value = "  +1.52%  "

$ws.Range("E23").Value = "  +0.82%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.23%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.18%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.62"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.12%  "

$ws.Range("E27").Value = "  +1.28%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.04%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0468"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.93%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.28%  "

$ws.Range("E32").Value = "  +0.73%  "

$ws.Range("D33").Value = "1.422.51"
$ws.Range("E33").Value = "  +0.25%  "

$ws.Range("E34").Value = "  +3.68%  "

$ws.Range("E35").Value = "  +3.55%  "

$ws.Range("E36").Value = "  +9.77%  "

$ws.Range("E38").Value = "  +0.75%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.532"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.92%  "

$ws.Range("E40").Value = "  +0.45%  "

$ws.Range("E41").Value = "  -0.12%  "

$ws.Range("E42").Value = "  +0.89%  "

$ws.Range("E43").Value = "  +2.97%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.996"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.75%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.01%  "

$ws.Range("E46").Value = "  +0.23%  "

$ws.Range("D47").Value = "1.696.61"
$ws.Range("E47").Value = "  +0.54%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.50%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0520"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.40%  "

$ws.Range("E50").Value = "  -0.01%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0956"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.30%  "
